# Fruta / hortaliza, semanal
# Update the weekly Espárragos records: Fecha (D) and Volumen (J) values are
# re-shuffled across rows 2-10, and the associated price columns (K/L/M/P)
# for the rows whose prices differed from the baseline move along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44166    # D2 (Fecha)
$ws.Cells.Item(2, 10).Value = 500     # J2 (Volumen)

# Row 3
$ws.Cells.Item(3, 4).Value = 44162    # D3
$ws.Cells.Item(3, 10).Value = 700     # J3

# Row 4
$ws.Cells.Item(4, 4).Value = 44165    # D4
$ws.Cells.Item(4, 10).Value = 300     # J4

# Row 5
$ws.Cells.Item(5, 4).Value = 44169    # D5
$ws.Cells.Item(5, 10).Value = 600     # J5

# Row 6
$ws.Cells.Item(6, 4).Value = 44159    # D6
$ws.Cells.Item(6, 12).Value = 1700    # L6 (Precio máximo)
$ws.Cells.Item(6, 13).Value = 1650    # M6 (Precio promedio ponderado)
$ws.Cells.Item(6, 16).Value = 1650    # P6 (Precio $/Kg)

# Row 7
$ws.Cells.Item(7, 4).Value = 44179    # D7
$ws.Cells.Item(7, 10).Value = 200     # J7
$ws.Cells.Item(7, 11).Value = 1600    # K7 (Precio mínimo)
$ws.Cells.Item(7, 12).Value = 1600    # L7
$ws.Cells.Item(7, 13).Value = 1600    # M7
$ws.Cells.Item(7, 16).Value = 1600    # P7

# Row 8
$ws.Cells.Item(8, 4).Value = 44168    # D8
$ws.Cells.Item(8, 10).Value = 200     # J8
$ws.Cells.Item(8, 12).Value = 1600    # L8
$ws.Cells.Item(8, 13).Value = 1600    # M8
$ws.Cells.Item(8, 16).Value = 1600    # P8

# Row 9
$ws.Cells.Item(9, 4).Value = 44176    # D9
$ws.Cells.Item(9, 10).Value = 700     # J9

# Row 10
$ws.Cells.Item(10, 4).Value = 44161   # D10
$ws.Cells.Item(10, 10).Value = 300    # J10
$ws.Cells.Item(10, 11).Value = 1700   # K10
$ws.Cells.Item(10, 12).Value = 1700   # L10
$ws.Cells.Item(10, 13).Value = 1700   # M10
$ws.Cells.Item(10, 16).Value = 1700   # P10
